# Fixed naive component forecaster bug - Presentation state 11.02.
#
# The data block B2:K16 (one row per forecast-origin date in column A, one
# column per forecast horizon Q0..Q9 in row 1) is missing the newly
# computed "first" horizon (Q0) error for each origin date. The fix
# inserts that new error into column B of every data row and pushes the
# previously-recorded errors one column to the right (the value that used
# to be the rightmost, Q9/column K, simply falls off the end of the
# table, since the table keeps only 10 horizon columns).
#
# Column A (the row date labels) and row 1 (the Q0..Q9 headers) are not
# touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing (pre-fix) horizon-error values for each row, B..K, as they
# appeared before this edit.
$oldRowValues = @{
    2  = @(-0.1220406607065741, -0.3824596977052969, 0.4529952984531944, 0.3529159385505728, -0.0008705675291234075, 0.1473068229033219, 0.3047528119074769, -0.1959066446185939, 0.08445123591687528, -0.08752362192813637)
    3  = @(-0.4000138190652613, 0.4354411770932299, 0.3353618171906084, -0.01842468888908786, 0.1297527015433575, 0.2871986905475125, -0.2134607659785583, 0.06689711455691082, -0.1050777432881008, 0.2007847964582307)
    4  = @(0.5235484960442361, 0.4234691361416146, 0.06968263006191837, 0.2178600204943637, 0.3753060094985187, -0.1253534470275521, 0.155004433507917, -0.01697042433709459, 0.2888921154092369, -0.3234378372780493)
    5  = @(0.6020437696293982, 0.248257263549702, 0.3964346539821473, 0.5538806429863024, 0.0532211864602315, 0.3335790669957007, 0.161604209150689, 0.4674667488970205, -0.1448632037902657, 0.5435081814088458)
    6  = @(1.160284491302352, 1.308461881734797, 1.465907870738952, 0.9652484142128814, 1.245606294748351, 1.073631436903339, 1.379493976649671, 0.7671640239623843, 1.455535409161496, 1.178844253737389)
    7  = @(0.3644876457690405, 0.5219336347731955, 0.02127417824712469, 0.3016320587825939, 0.1296572009375822, 0.4355197406839137, -0.1768102120033725, 0.511561173195739, 0.2348700177716323)
    8  = @(0.5259015322863045, 0.0252420757602338, 0.305599956295703, 0.1336250984506913, 0.4394876381970228, -0.1728423144902634, 0.5155290707088481, 0.2388379152847414)
    9  = @(0.1608821659304752, 0.4412400464659443, 0.2692651886209327, 0.5751277283672642, -0.03720222432002201, 0.6511691608790895, 0.3744780054549828)
    10 = @(0.2004338646102796, 0.028459006765268, 0.3343215465115995, -0.2780084061756867, 0.4103629790234248, 0.1336718235993181)
    11 = @(-0.01687220848682837, 0.2889903312595031, -0.3233396214277831, 0.3650317637713285, 0.08834060834722172)
    12 = @(0.2221289093234493, -0.3902010433638369, 0.2981703418352746, 0.02147918641116785)
    13 = @(-0.4197872457205535, 0.268584139478558, -0.00810701594554874)
    14 = @(0.2504334827489171, -0.02625767267518964)
    15 = @(-0.04428949692388896)
    16 = @()
}

# Newly computed first-horizon (Q0 / column B) error for each row.
$newFirstValue = @{
    2  = -0.2177157015159319
    3  = -0.1395947820665385
    4  = -0.3119065001142551
    5  = 0.7021231295320197
    6  = 1.514070997382048
    7  = 0.2163102553365951
    8  = 0.3684555432821496
    9  = 0.661541622456546
    10 = -0.07992401592518952
    11 = 0.1551026493581833
    12 = -0.08373363042288225
    13 = 0.1925427069667326
    14 = -0.4379379024501944
    15 = 0.2324016585002178
    16 = -0.09587373626955231
}

$columns = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

foreach ($row in 2..16) {
    # Build the post-fix row: new value first, then the old values shifted
    # right, truncated to the 10 available horizon columns (B..K).
    $rowValues = @($newFirstValue[$row]) + $oldRowValues[$row]
    if ($rowValues.Count -gt $columns.Count) {
        $rowValues = $rowValues[0..($columns.Count - 1)]
    }

    for ($i = 0; $i -lt $rowValues.Count; $i++) {
        $ws.Range($columns[$i] + $row).Value = $rowValues[$i]
    }
}
